$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 4576.875
$ws.Range("I15").Value = 4576.875
$ws.Range("K15").Value = 13730.625
$ws.Range("M15").Value = -13561.625
$ws.Range("H32").Value = 6839.8
$ws.Range("I32").Value = 6839.8
$ws.Range("K32").Value = 6839.8
$ws.Range("M32").Value = -6513.8
$ws.Range("H40").Value = 4610.35
$ws.Range("I40").Value = 3067.5
$ws.Range("J40").Value = 4996.0625
$ws.Range("K40").Value = 3067.5
$ws.Range("L40").Value = 4996.0625
$ws.Range("M40").Value = -2892.5
$ws.Range("N40").Value = -5346.0625
$ws.Range("H45").Value = 225.8
$ws.Range("J45").Value = 225.8
$ws.Range("L45").Value = 677.4000000000001
$ws.Range("N45").Value = -1061.4
$ws.Range("H49").Value = 740.8333
$ws.Range("J49").Value = 740.8333
$ws.Range("L49").Value = 2222.4999
$ws.Range("N49").Value = -2494.4999
$ws.Range("H51").Value = 4825.5557
$ws.Range("I51").Value = 5113.75
$ws.Range("J51").Value = 4595
$ws.Range("K51").Value = 5113.75
$ws.Range("L51").Value = 4595
$ws.Range("M51").Value = -4629.75
$ws.Range("N51").Value = -5563
$ws.Range("H53").Value = 111838.39
$ws.Range("I53").Value = 926.75
$ws.Range("J53").Value = 200567.7
$ws.Range("K53").Value = 926.75
$ws.Range("L53").Value = 200567.7
$ws.Range("M53").Value = -289.75
$ws.Range("N53").Value = -201841.7
$ws.Range("H98").Value = 761.63635
$ws.Range("J98").Value = 733.8
$ws.Range("L98").Value = 733.8
$ws.Range("N98").Value = -3729.8
$ws.Range("H100").Value = 1332.2858
$ws.Range("I100").Value = 1332.2858
$ws.Range("K100").Value = 1332.2858
$ws.Range("M100").Value = -791.2858000000001
$ws.Range("H107").Value = 742.3125
$ws.Range("J107").Value = 983
$ws.Range("L107").Value = 983
$ws.Range("N107").Value = -4823
$ws.Range("H112").Value = 3684.75
$ws.Range("J112").Value = 4299.8
$ws.Range("L112").Value = 12899.4
$ws.Range("N112").Value = -15115.4
$ws.Range("H122").Value = 761.63635
$ws.Range("J122").Value = 733.8
$ws.Range("L122").Value = 2201.4
$ws.Range("N122").Value = -7101.4
$ws.Range("H131").Value = 5636.125
$ws.Range("I131").Value = 2348.1667
$ws.Range("K131").Value = 7044.500100000001
$ws.Range("M131").Value = -2004.500100000001
$ws.Range("H132").Value = 1543.5366
$ws.Range("I132").Value = 1541.3611
$ws.Range("K132").Value = 4624.0833
$ws.Range("M132").Value = -2094.0833

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5168.2075
$ws.Range("I32").Value = 4418.04
$ws.Range("K32").Value = 4418.04
$ws.Range("M32").Value = -4131.04
$ws.Range("H74").Value = 6835.5625
$ws.Range("I74").Value = 4124.5454
$ws.Range("K74").Value = 4124.5454
$ws.Range("M74").Value = -3250.5454
$ws.Range("H77").Value = 6835.5625
$ws.Range("I77").Value = 4124.5454
$ws.Range("K77").Value = 20622.727
$ws.Range("M77").Value = -16254.727
$ws.Range("H135").Value = 68929.11
$ws.Range("J135").Value = 68929.11
$ws.Range("L135").Value = 68929.11
$ws.Range("N135").Value = -79069.11

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 69980
$ws.Range("J2").Value = 69980
$ws.Range("L2").Value = 69980
$ws.Range("N2").Value = -70206
$ws.Range("H26").Value = 20471
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5282.1934
$ws.Range("I31").Value = 3846
$ws.Range("K31").Value = 3846
$ws.Range("M31").Value = -3551
$ws.Range("H34").Value = 5282.1934
$ws.Range("I34").Value = 3846
$ws.Range("K34").Value = 3846
$ws.Range("M34").Value = -3644
$ws.Range("H58").Value = 12328.857
$ws.Range("I58").Value = 8163.3335
$ws.Range("J58").Value = 13464.909
$ws.Range("K58").Value = 8163.3335
$ws.Range("L58").Value = 13464.909
$ws.Range("M58").Value = -7960.3335
$ws.Range("N58").Value = -13870.909
$ws.Range("H64").Value = 99999
$ws.Range("J64").Value = 99999
$ws.Range("L64").Value = 99999
$ws.Range("N64").Value = -100495
$ws.Range("H67").Value = 99999
$ws.Range("J67").Value = 99999
$ws.Range("L67").Value = 99999
$ws.Range("N67").Value = -101715
$ws.Range("H105").Value = 1141
$ws.Range("I105").Value = 1311.75
$ws.Range("J105").Value = 799.5
$ws.Range("K105").Value = 1311.75
$ws.Range("L105").Value = 799.5
$ws.Range("M105").Value = 435.25
$ws.Range("N105").Value = -4293.5
$ws.Range("H125").Value = 36445
$ws.Range("J125").Value = 36445
$ws.Range("L125").Value = 36445
$ws.Range("N125").Value = -41365
$ws.Range("H132").Value = 8350
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 8350
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 25050
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -30110
$ws.Range("H136").Value = 12328.857
$ws.Range("I136").Value = 8163.3335
$ws.Range("J136").Value = 13464.909
$ws.Range("K136").Value = 24490.0005
$ws.Range("L136").Value = 40394.727
$ws.Range("M136").Value = -21940.0005
$ws.Range("N136").Value = -45494.727
$ws.Range("H141").Value = 316967.84
$ws.Range("J141").Value = 316967.84
$ws.Range("L141").Value = 316967.84
$ws.Range("N141").Value = -327327.84

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 354.22223
$ws.Range("I25").Value = 99.5
$ws.Range("J25").Value = 427
$ws.Range("K25").Value = 298.5
$ws.Range("L25").Value = 1281
$ws.Range("M25").Value = -129.5
$ws.Range("N25").Value = -1619
$ws.Range("H30").Value = 354.22223
$ws.Range("I30").Value = 99.5
$ws.Range("J30").Value = 427
$ws.Range("K30").Value = 298.5
$ws.Range("L30").Value = 1281
$ws.Range("M30").Value = -196.5
$ws.Range("N30").Value = -1485
$ws.Range("H40").Value = 152.06667
$ws.Range("I40").Value = 96.75
$ws.Range("K40").Value = 387
$ws.Range("M40").Value = -318
$ws.Range("H61").Value = 164.625
$ws.Range("I61").Value = 78
$ws.Range("K61").Value = 234
$ws.Range("M61").Value = -19

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 33999
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 33999
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 33999
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -35135
$ws.Range("H62").Value = 41949
$ws.Range("I62").Value = 41949
$ws.Range("K62").Value = 41949
$ws.Range("M62").Value = -41263
$ws.Range("H65").Value = 41949
$ws.Range("I65").Value = 41949
$ws.Range("K65").Value = 125847
$ws.Range("M65").Value = -122415
$ws.Range("H123").Value = 54612.5
$ws.Range("J123").Value = 54612.5
$ws.Range("L123").Value = 54612.5
$ws.Range("N123").Value = -59512.5
$ws.Range("H132").Value = 5706.2856
$ws.Range("I132").Value = 5534.1816
$ws.Range("J132").Value = 6337.3335
$ws.Range("K132").Value = 16602.5448
$ws.Range("L132").Value = 19012.0005
$ws.Range("M132").Value = -14072.5448
$ws.Range("N132").Value = -24072.0005
$ws.Range("H136").Value = 38698.438
$ws.Range("J136").Value = 38698.438
$ws.Range("L136").Value = 116095.314
$ws.Range("N136").Value = -121195.314

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2681.9
$ws.Range("J46").Value = 3899.8333
$ws.Range("L46").Value = 3899.8333
$ws.Range("N46").Value = -4275.8333
$ws.Range("H82").Value = 1250.5
$ws.Range("I82").Value = 1391.5714
$ws.Range("K82").Value = 1391.5714
$ws.Range("M82").Value = -1030.5714
$ws.Range("H85").Value = 1250.5
$ws.Range("I85").Value = 1391.5714
$ws.Range("K85").Value = 1391.5714
$ws.Range("M85").Value = -143.5714
$ws.Range("H100").Value = 5249.625
$ws.Range("H132").Value = 7730.316
$ws.Range("J132").Value = 3929.2
$ws.Range("L132").Value = 11787.6
$ws.Range("N132").Value = -16847.6
$ws.Range("H136").Value = 6976.0586
$ws.Range("I136").Value = 6976.0586
$ws.Range("K136").Value = 20928.1758
$ws.Range("M136").Value = -18378.1758

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H76").Value = 68000
$ws.Range("J76").Value = 68000
$ws.Range("L76").Value = 68000
$ws.Range("N76").Value = -68630
$ws.Range("H79").Value = 68000
$ws.Range("J79").Value = 68000
$ws.Range("L79").Value = 68000
$ws.Range("N79").Value = -70184
$ws.Range("H132").Value = 5510.24
$ws.Range("I132").Value = 4625.263
$ws.Range("K132").Value = 13875.789
$ws.Range("M132").Value = -11345.789
$ws.Range("H136").Value = 2914.182
$ws.Range("I136").Value = 2814.639
$ws.Range("K136").Value = 8443.917000000001
$ws.Range("M136").Value = -5893.917000000001
